$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.655.90"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").Value = "1.620.10"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.992"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.84%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.73"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.516"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.42%  "
$ws.Range("E7").Value = "  -0.88%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.19"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.66%  "
$ws.Range("E9").Value = "  -1.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0606"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.45%  "
$ws.Range("E11").Value = "  -0.99%  "
$ws.Range("D12").Value = "1.852.08"
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("D13").Value = "1.622.58"
$ws.Range("E13").Value = "  -0.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.98"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.557"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.74"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("D17").Value = "27.690.97"
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "227.83"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.37%  "
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("D20").Value = "0.0₃0715"
$ws.Range("E20").Value = "  -1.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.993"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.31"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.07"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.09%  "
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.51"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.90"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.94%  "
$ws.Range("E27").Value = "  -0.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.43"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.993"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("E30").Value = "  -1.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0477"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.38"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.73%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").Value = "1.390.48"
$ws.Range("E34").Value = "  -1.18%  "
$ws.Range("E35").Value = "  +1.50%  "
$ws.Range("E36").Value = "  -0.79%  "
$ws.Range("E37").Value = "  -1.45%  "
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.556"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.844"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.27%  "
$ws.Range("E41").Value = "  -1.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.992"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.82"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "65.56"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.99%  "
$ws.Range("E45").Value = "  -2.74%  "
$ws.Range("D46").Value = "1.761.18"
$ws.Range("E46").Value = "  -0.92%  "
$ws.Range("E47").Value = "  -1.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.74"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.36%  "
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0502"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.58"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.15%  "
